$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column F (dSF) values for the rows that changed,
# reflecting a repulled dataset / recalculated mean.
$ws.Range("F3").Value = -1
$ws.Range("F5").Value = -3
$ws.Range("F6").Value = -2
$ws.Range("F7").Value = 5
$ws.Range("F8").Value = -6
$ws.Range("F10").Value = -7
$ws.Range("F11").Value = -1
$ws.Range("F12").Value = 0
